$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{}
$data[2] = @('Insight=Religion', 'Insight=Religion', 'History=Survival', 'Acrobatics=Athletics', 'Insight=Intimidation', 'Deception=Insight', 'Athletics=Intimidation', 'Animal Handling=Survival', 'Perception=Performance', 'Deception=Sleight of Hand', 'Athletics=Insight')
$data[3] = @('None', 'None', 'Navigator’s tools', 'Vehicles (Land)', 'None', 'Disguise Kit=One Type of Artisan’s Tools', 'One Type of Gaming Set', 'Land Vehicles', 'Disguise Kit', 'Disguise Kit=Forgery Kit', 'None')
$data[4] = @('Any=Any', 'Any=Any', 'Any', 'Any', 'Any=Any', 'None', 'Celestial=Draconic', 'Any', 'Any=Any', 'None', 'Any=Any')
$data[5] = @('A Holy Symbol=A Prayer Book=5 Sticks of Incense=Vestments=a Set of Common Clothes=A Belt Pouch which is Containing 15 gp', 'A Leather-Bound Diary=A Bottle of Ink=An Ink Pen=A Set of Traveler’s Clothes=One Trinket of Special Significance=A Pouch Containing 10 gp', 'A Wooden Case Containing a Map to a Ruin=A Bullseye Lantern=A Miner’s Pick=A Set of Traveler’s Clothes=A Shovel=A Two-Person Tent=A Trinket Recovered from a Dig Site=A Pouch Containing 25 gp', 'A Leather Ball=A Lucky Charm=A Set of Traveler’s Clothes=A Pouch Containing 10 gp', 'An Azorius Insignia=A Scroll Containing the Text of a Law Important to you=A Bottle of Blue Ink=A Pen=A Set of Fine Clothes=A Belt Pouch Containing 10 gp (Azorius-minted 1-zino coins)', 'Disguise Kit=Common Clothes=A Tears of Virulence Emblem=A Writ of Free Agency Signed by the Lord Regent=A Set of Artisan’s Tools you are Proficient With=A Pouch with 15 gp (payment for services rendered)', 'A Boros Insignia=A Feather from an Angel’s Wing=A Tattered Piece of a Boros Banner (A Souvenir from a Famous Battle)=A Set of Common Clothes=A Belt Pouch Containing 2gp (Boros-minted 1-zino coins)', 'A Whip=A Tent=A Regional Map=A Set of Traveling Clothes=Belt Pouch Containing 10 gp', 'Disguise kit=A Set of Fine Clothes=Belt Pouch Containing 30 gp(gold piece)', 'A Set of Fine Clothes=Disguise Kit=Tools of the Con of your Choice (Ten Stoppered Bottles filled with Colored Liquid, a Set of Weighted Dice, a Deck of Marked Cards, or a Signet Ring of an Imaginary Duke)=Belt Pouch Containing 15 gp', 'A Uniform in the Style of your Unit=Indicative of your Rank=Horn with which to Summon Help=A Set of Manacles=A Pouch Containing 10 gp')

foreach ($r in 2..5) {
    $vals = $data[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $cell = $ws.Cells.Item($r, $i + 1)
        $cell.Value = $vals[$i]
    }
}

# Apply the new cell style (font: Arial, theme color 1) used for the new data rows
$ws.Range("A2:K5").Font.ThemeColor = 1
